# Applies:
#   1) Refresh the "panel_query_time"-style timestamps in column F of the
#      "data" sheet (rows 2-48) to the new query run's values.
#   2) Add a new "metadata" worksheet after "data" describing the panel
#      query itself (name/id/version/etc.), mirroring the header style used
#      on the "data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Update time_taken values on the "data" sheet (F2:F48) -------------
$times = @(
  "2021-10-05 14:20:08.814786",
  "2021-10-05 14:20:08.814795",
  "2021-10-05 14:20:08.814799",
  "2021-10-05 14:20:08.814802",
  "2021-10-05 14:20:08.814804",
  "2021-10-05 14:20:08.814807",
  "2021-10-05 14:20:08.814810",
  "2021-10-05 14:20:08.814812",
  "2021-10-05 14:20:08.814815",
  "2021-10-05 14:20:08.814818",
  "2021-10-05 14:20:08.814820",
  "2021-10-05 14:20:08.814823",
  "2021-10-05 14:20:08.814825",
  "2021-10-05 14:20:08.814828",
  "2021-10-05 14:20:08.814830",
  "2021-10-05 14:20:08.814833",
  "2021-10-05 14:20:08.814835",
  "2021-10-05 14:20:08.814838",
  "2021-10-05 14:20:08.814840",
  "2021-10-05 14:20:08.814843",
  "2021-10-05 14:20:08.814845",
  "2021-10-05 14:20:08.814848",
  "2021-10-05 14:20:08.814850",
  "2021-10-05 14:20:08.814853",
  "2021-10-05 14:20:08.814856",
  "2021-10-05 14:20:08.814858",
  "2021-10-05 14:20:08.814861",
  "2021-10-05 14:20:08.814863",
  "2021-10-05 14:20:08.814866",
  "2021-10-05 14:20:08.814868",
  "2021-10-05 14:20:08.814871",
  "2021-10-05 14:20:08.814873",
  "2021-10-05 14:20:08.814876",
  "2021-10-05 14:20:08.814879",
  "2021-10-05 14:20:08.814881",
  "2021-10-05 14:20:08.814884",
  "2021-10-05 14:20:08.814886",
  "2021-10-05 14:20:08.814889",
  "2021-10-05 14:20:08.814891",
  "2021-10-05 14:20:08.814893",
  "2021-10-05 14:20:08.814896",
  "2021-10-05 14:20:08.814899",
  "2021-10-05 14:20:08.814902",
  "2021-10-05 14:20:08.814904",
  "2021-10-05 14:20:08.814906",
  "2021-10-05 14:20:08.814909",
  "2021-10-05 14:20:08.814911"
)

for ($i = 0; $i -lt $times.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}

# --- 2) Add the "metadata" sheet, placed right after "data" ---------------
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (columns B:G), same layout style as the "data" sheet (column A
# is reserved for the numeric index column).
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Single metadata record (row 2).
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Epidermolysis bullosa and congenital skin fragility"
$meta.Range("C2").Value = 554
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.49"
$meta.Range("D2").Style = "Normal"
$meta.Range("E2").Value = "2021-07-19T09:09:24.950869Z"
$meta.Range("F2").Value = "2021-10-05 14:20:08.810957"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/554/?format=json"

# Match the bold / centered / bordered header formatting already used by the
# "data" sheet's header row (B1) and its index column (A2).
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep "data" selected/active, matching the unmodified bookView state.
$ws.Activate()
